$wb = $excel.ActiveWorkbook

# Sheet "DatosCuenta"
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneLastTwo"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneLastTwo"
$wsCuenta.Range("C2").Value = 27100133
$wsCuenta.Range("D2").Value = 134

# Sheet "DatosHogar"
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 653

# Sheet "DatosMotor"
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP035"
$wsMotor.Range("B2").Value = "ABC12SSMP035"
$wsMotor.Range("C2").Value = "ZAZ123SSMP035"

# Sheet "DatosAP"
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200134
$wsAP.Activate()
$wsAP.Range("E6").Select()
